$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.928.49"
$ws.Range("E2").Value = "  +0.64%  "

$ws.Range("D3").Value = "1.924.85"
$ws.Range("E3").Value = "  +1.81%  "

$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").Value = "240.66"
$ws.Range("E5").Value = "  -2.92%  "

$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.15%  "

$ws.Range("D7").Value = "0.4911"
$ws.Range("E7").Value = "  -0.71%  "

$ws.Range("D8").Value = "0.2983"
$ws.Range("E8").Value = "  +0.62%  "

$ws.Range("D9").Value = "0.06781"
$ws.Range("E9").Value = "  -0.73%  "

$ws.Range("D10").Value = "1.915.85"
$ws.Range("E10").Value = "  +1.34%  "

$ws.Range("D11").Value = "17.11"

$ws.Range("D12").Value = "0.07302"
$ws.Range("E12").Value = "  +0.84%  "

$ws.Range("D13").Value = "5.178"
$ws.Range("E13").Value = "  +1.92%  "

$ws.Range("D14").Value = "89.71"
$ws.Range("E14").Value = "  -2.56%  "

$ws.Range("D15").Value = "0.6734"
$ws.Range("E15").Value = "  -0.88%  "

$ws.Range("D16").Value = "30.898.90"
$ws.Range("E16").Value = "  +0.63%  "

$ws.Range("D17").Value = "0.000007998"
$ws.Range("E17").Value = "  +0.22%  "

$ws.Range("E18").Value = "  +2.46%  "

$ws.Range("D19").Value = "1.002"
$ws.Range("E19").Value = "  +0.16%  "

$ws.Range("D20").Value = "2.179.52"
$ws.Range("E20").Value = "  +2.06%  "

$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  +0.15%  "

$ws.Range("D22").Value = "5.178"
$ws.Range("E22").Value = "  +6.93%  "

$ws.Range("D23").Value = "206.64"
$ws.Range("E23").Value = "  +8.16%  "

$ws.Range("D24").Value = "6.326"
$ws.Range("E24").Value = "  +4.24%  "

$ws.Range("D25").Value = "9.706"

$ws.Range("D26").Value = "159.30"
$ws.Range("E26").Value = "  +1.52%  "

$ws.Range("D27").Value = "19.11"
$ws.Range("E27").Value = "  +0.22%  "

$ws.Range("D28").Value = "1.996"
$ws.Range("E28").Value = "  +4.20%  "

$ws.Range("D29").Value = "1.429"
$ws.Range("E29").Value = "  +1.82%  "

$ws.Range("D30").Value = "4.378"
$ws.Range("E30").Value = "  +1.55%  "

$ws.Range("D31").Value = "0.09203"
$ws.Range("E31").Value = "  +2.06%  "

$ws.Range("D32").Value = "4.085"
$ws.Range("E32").Value = "  +1.63%  "

$ws.Range("D33").Value = "0.05209"
$ws.Range("E33").Value = "  +0.28%  "

$ws.Range("D34").Value = "0.7550"
$ws.Range("E34").Value = "  +1.28%  "

$ws.Range("D35").Value = "1.126"
$ws.Range("E35").Value = "  +0.56%  "

$ws.Range("D36").Value = "2.731"
$ws.Range("E36").Value = "  +0.31%  "

$ws.Range("D37").Value = "0.01863"

$ws.Range("D38").Value = "2.741"
$ws.Range("E38").Value = "  +2.47%  "

$ws.Range("D39").Value = "0.9293"
$ws.Range("E39").Value = "  -1.34%  "

$ws.Range("D40").Value = "2.099"
$ws.Range("E40").Value = "  -3.22%  "

$ws.Range("D41").Value = "0.4532"
$ws.Range("E41").Value = "  +2.22%  "

$ws.Range("E42").Value = "  +2.39%  "

$ws.Range("D43").Value = "5.953"
$ws.Range("E43").Value = "  +3.10%  "

$ws.Range("D44").Value = "71.93"
$ws.Range("E44").Value = "  +23.93%  "

$ws.Range("E45").Value = "  +1.17%  "

$ws.Range("D46").Value = "0.1398"
$ws.Range("E46").Value = "  +3.99%  "

$ws.Range("D47").Value = "7.703"
$ws.Range("E47").Value = "  +0.57%  "

$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").Value = "35.66"
$ws.Range("E48").Value = "  +6.27%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "9.107"
$ws.Range("E49").Value = "  +4.76%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.05961"
$ws.Range("E50").Value = "  +2.04%  "

$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").Value = "0.4099"
$ws.Range("E51").Value = "  +3.71%  "
